$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.966.52'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '1.593.26'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  -0.30%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '210.51'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("E7").Value = '  +0.27%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.246'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.32%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.0611'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.89%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '17.96'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.74%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0809'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +2.39%  '
$ws.Range("D12").Value = '1.815.84'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = '1.592.57'
$ws.Range("E13").Value = '  +0.35%  '
$ws.Range("E14").Value = '  -0.93%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.513'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("D16").Value = '25.982.64'
$ws.Range("E16").Value = '  +0.65%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '60.05'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("E19").Value = '  -0.21%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '199.57'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +3.80%  '
$ws.Range("E21").Value = '  +0.84%  '
$ws.Range("E22").Value = '  -1.47%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.99'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.19%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '1.86'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +9.64%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '143.25'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("E27").Value = '  -7.50%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '15.09'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.30%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("E32").Value = '  +0.15%  '
$ws.Range("E33").Value = '  -2.51%  '
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("D36").Value = '1.122.97'
$ws.Range("E36").Value = '  +2.03%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.0162'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +8.27%  '
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("E39").Value = '  -1.34%  '
$ws.Range("E40").Value = '  +0.62%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.490'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -2.69%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.781'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -4.02%  '
$ws.Range("D43").Value = '1.727.21'
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("E44").Value = '  -1.45%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '92.50'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.36%  '
$ws.Range("E46").Value = '  -0.84%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '53.32'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("E48").Value = '  -1.06%  '
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("D51").Value = '0.0₇0916'
$ws.Range("E51").Value = '  -14.06%  '
